$wb = $excel.ActiveWorkbook

# Helper: write a numeric-looking value as TEXT (shared string), preserving the
# cell's existing style (avoids Excel's "looks like a number" auto-conversion,
# and avoids adding a quote-prefix style that a plain apostrophe would add).
function Set-TextValue($range, $text) {
    $range.Formula = '="' + $text + '"'
    $range.Copy()
    $range.PasteSpecial(-4163, $null)   # xlPasteValues
}

# Sheet1 ("Sheet1" tab) - module/date/mobile columns regenerated with new test data
$ws1 = $wb.Worksheets.Item(1)
Set-TextValue $ws1.Range("F2")  "2007632789"
Set-TextValue $ws1.Range("N2")  "2024-02-16"
Set-TextValue $ws1.Range("P2")  "2024-02-16 03:04:56 PM"
Set-TextValue $ws1.Range("AC2") "2024-02-16"
Set-TextValue $ws1.Range("AE2") "0618268283"
Set-TextValue $ws1.Range("AN2") "94424"
Set-TextValue $ws1.Range("AT2") "7630880617"
Set-TextValue $ws1.Range("AX2") "8581647506"

# Sheet2 - mobile/phone columns regenerated with new test data
$ws2 = $wb.Worksheets.Item(2)
Set-TextValue $ws2.Range("F2")  "2007632789"
Set-TextValue $ws2.Range("AE2") "0618268283"
Set-TextValue $ws2.Range("AT2") "7630880617"
Set-TextValue $ws2.Range("AX2") "8581647506"

# Sheet3 - mobile/phone columns regenerated with new test data
$ws3 = $wb.Worksheets.Item(3)
Set-TextValue $ws3.Range("F2")  "2007632789"
Set-TextValue $ws3.Range("AE2") "0618268283"
Set-TextValue $ws3.Range("AT2") "7630880617"
Set-TextValue $ws3.Range("AX2") "8581647506"

# Sheet4 - mobile/phone columns regenerated with new test data
$ws4 = $wb.Worksheets.Item(4)
Set-TextValue $ws4.Range("F2")  "2007632789"
Set-TextValue $ws4.Range("AE2") "0618268283"
Set-TextValue $ws4.Range("AT2") "7630880617"
Set-TextValue $ws4.Range("AX2") "8581647506"
